# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Femacal de La Calera - Pepino ensalada"
# right after the existing row 181 (i.e. at row 182), shifting all the
# following rows down by one. The sheet's used range grows from
# A1:R227 to A1:R228.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 182:227 down to 183:228 by inserting a blank row at 182.
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(182, 1).Value = 3
$ws.Cells.Item(182, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44508
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 100112043
$ws.Cells.Item(182, 7).Value = "Pepino ensalada"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 85
$ws.Cells.Item(182, 11).Value = 7000
$ws.Cells.Item(182, 12).Value = 7500
$ws.Cells.Item(182, 13).Value = 7235
$ws.Cells.Item(182, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(182, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(182, 16).Value = 103
$ws.Cells.Item(182, 17).Value = 70
$ws.Cells.Item(182, 18).Value = "Hortaliza"
